$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.550.02'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '3.362.84'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '256.50'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.21%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '663.20'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +5.95%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.54'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +4.10%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.473'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +19.63%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.08'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +18.92%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('D11').Value = '3.360.19'
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.218'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +8.76%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '42.53'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +9.19%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000273'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +9.30%  '
$ws.Range('D15').Value = '98.866.24'
$ws.Range('E15').Value = '  +0.35%  '
$ws.Range('B16').Value = 'Toncoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.75'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.69%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '3.996.97'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.97'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +28.82%  '
$ws.Range('D19').Value = '3.358.24'
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.14'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +10.98%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '530.06'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +8.27%  '
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.49'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +8.62%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.0000216'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.72%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.443'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +47.71%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '103.29'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +15.23%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.30'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +11.61%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '12.66'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.78%  '
$ws.Range('D29').Value = '3.537.53'
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.154'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +12.86%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '11.42'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +15.42%  '
$ws.Range('B32').Value = 'Dai'
$ws.Range('C32').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.999'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.194'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '29.65'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +4.03%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.556'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +18.19%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.13'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +8.95%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '7.72'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +4.81%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.158'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +5.30%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '524.24'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +4.09%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.34'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +5.36%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '24.74'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('B43').Value = 'MantraDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.85'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.29%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0437'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +31.09%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.44'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.25%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.834'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.09'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +6.83%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.03'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +15.67%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '5.23'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +10.88%  '
$ws.Range('B51').Value = 'ImmutableX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.55'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +12.23%  '
